$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.678587
$ws.Range("H2").Value = 65.03576100000001
$ws.Range("I2").Value = 0.0772399443186744
$ws.Range("J2").Value = 0.07723994431867441
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.222961333333333
$ws.Range("N2").Value = 3.668884
$ws.Range("O2").Value = 0.02082890820948692
$ws.Range("P2").Value = 0.02082890820948692
$ws.Range("Q2").Value = 26.51207366230267
$ws.Range("R2").Value = 238.608662960724
$ws.Range("S2").Value = 0.00160882371031955
$ws.Range("T2").Value = 0.00160882371031955
$ws.Range("G3").Value = 21.678587
$ws.Range("H3").Value = 65.03576100000001
$ws.Range("I3").Value = 0.0772399443186744
$ws.Range("J3").Value = 0.07723994431867441
$ws.Range("O3").Value = 0.1691870972318839
$ws.Range("P3").Value = 0.169187097231884
$ws.Range("Q3").Value = 215.3497792303807
$ws.Range("R3").Value = 1938.148013073426
$ws.Range("S3").Value = 0.01306800196962887
$ws.Range("T3").Value = 0.01306800196962887
$ws.Range("G4").Value = 21.678587
$ws.Range("H4").Value = 65.03576100000001
$ws.Range("I4").Value = 0.0772399443186744
$ws.Range("J4").Value = 0.07723994431867441
$ws.Range("M4").Value = 14.516908
$ws.Range("N4").Value = 43.550724
$ws.Range("O4").Value = 0.247245220250272
$ws.Range("P4").Value = 0.2472452202502721
$ws.Range("Q4").Value = 314.7060530489961
$ws.Range("R4").Value = 2832.354477440965
$ws.Range("S4").Value = 0.0190972070451894
$ws.Range("T4").Value = 0.01909720704518941
$ws.Range("G5").Value = 21.678587
$ws.Range("H5").Value = 65.03576100000001
$ws.Range("I5").Value = 0.0772399443186744
$ws.Range("J5").Value = 0.07723994431867441
$ws.Range("M5").Value = 12.24131666666667
$ws.Range("N5").Value = 36.72395
$ws.Range("O5").Value = 0.2084884078209579
$ws.Range("P5").Value = 0.2084884078209579
$ws.Range("Q5").Value = 265.3744483528834
$ws.Range("R5").Value = 2388.37003517595
$ws.Range("S5").Value = 0.01610363301117987
$ws.Range("T5").Value = 0.01610363301117987
$ws.Range("G6").Value = 21.678587
$ws.Range("H6").Value = 65.03576100000001
$ws.Range("I6").Value = 0.0772399443186744
$ws.Range("J6").Value = 0.07723994431867441
$ws.Range("M6").Value = 20.799674
$ws.Range("N6").Value = 62.399022
$ws.Range("O6").Value = 0.3542503664873991
$ws.Range("P6").Value = 0.3542503664873992
$ws.Range("Q6").Value = 450.9075423806381
$ws.Range("R6").Value = 4058.167881425743
$ws.Range("S6").Value = 0.0273622785823567
$ws.Range("T6").Value = 0.02736227858235672
$ws.Range("I7").Value = 0.4677505770609061
$ws.Range("J7").Value = 0.4677505770609062
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.222961333333333
$ws.Range("N7").Value = 3.668884
$ws.Range("O7").Value = 0.02082890820948692
$ws.Range("P7").Value = 0.02082890820948692
$ws.Range("Q7").Value = 160.552132242088
$ws.Range("R7").Value = 1444.969190178792
$ws.Range("S7").Value = 0.009742733834536148
$ws.Range("T7").Value = 0.009742733834536155
$ws.Range("I8").Value = 0.4677505770609061
$ws.Range("J8").Value = 0.4677505770609062
$ws.Range("O8").Value = 0.1691870972318839
$ws.Range("P8").Value = 0.169187097231884
$ws.Range("S8").Value = 0.07913736236147334
$ws.Range("T8").Value = 0.07913736236147338
$ws.Range("I9").Value = 0.4677505770609061
$ws.Range("J9").Value = 0.4677505770609062
$ws.Range("M9").Value = 14.516908
$ws.Range("N9").Value = 43.550724
$ws.Range("O9").Value = 0.247245220250272
$ws.Range("P9").Value = 0.2472452202502721
$ws.Range("Q9").Value = 1905.800673688968
$ws.Range("R9").Value = 17152.20606320071
$ws.Range("S9").Value = 0.1156490944476155
$ws.Range("T9").Value = 0.1156490944476156
$ws.Range("I10").Value = 0.4677505770609061
$ws.Range("J10").Value = 0.4677505770609062
$ws.Range("M10").Value = 12.24131666666667
$ws.Range("N10").Value = 36.72395
$ws.Range("O10").Value = 0.2084884078209579
$ws.Range("P10").Value = 0.2084884078209579
$ws.Range("Q10").Value = 1607.0577529439
$ws.Range("R10").Value = 14463.5197764951
$ws.Range("S10").Value = 0.09752057306876256
$ws.Range("T10").Value = 0.09752057306876261
$ws.Range("I11").Value = 0.4677505770609061
$ws.Range("J11").Value = 0.4677505770609062
$ws.Range("M11").Value = 20.799674
$ws.Range("N11").Value = 62.399022
$ws.Range("O11").Value = 0.3542503664873991
$ws.Range("P11").Value = 0.3542503664873992
$ws.Range("Q11").Value = 2730.611279048604
$ws.Range("R11").Value = 24575.50151143744
$ws.Range("S11").Value = 0.1657008133485184
$ws.Range("T11").Value = 0.1657008133485185
$ws.Range("G12").Value = 67.370907
$ws.Range("H12").Value = 202.112721
$ws.Range("I12").Value = 0.2400398653924534
$ws.Range("J12").Value = 0.2400398653924535
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 1.222961333333333
$ws.Range("N12").Value = 3.668884
$ws.Range("O12").Value = 0.02082890820948692
$ws.Range("P12").Value = 0.02082890820948692
$ws.Range("Q12").Value = 82.392014252596
$ws.Range("R12").Value = 741.5281282733641
$ws.Range("S12").Value = 0.004999768322877007
$ws.Range("T12").Value = 0.004999768322877011
$ws.Range("G13").Value = 67.370907
$ws.Range("H13").Value = 202.112721
$ws.Range("I13").Value = 0.2400398653924534
$ws.Range("J13").Value = 0.2400398653924535
$ws.Range("O13").Value = 0.1691870972318839
$ws.Range("P13").Value = 0.169187097231884
$ws.Range("Q13").Value = 669.2461067227539
$ws.Range("R13").Value = 6023.214960504786
$ws.Range("S13").Value = 0.04061164804568135
$ws.Range("T13").Value = 0.04061164804568137
$ws.Range("G14").Value = 67.370907
$ws.Range("H14").Value = 202.112721
$ws.Range("I14").Value = 0.2400398653924534
$ws.Range("J14").Value = 0.2400398653924535
$ws.Range("M14").Value = 14.516908
$ws.Range("N14").Value = 43.550724
$ws.Range("O14").Value = 0.247245220250272
$ws.Range("P14").Value = 0.2472452202502721
$ws.Range("Q14").Value = 978.0172587955561
$ws.Range("R14").Value = 8802.155329160005
$ws.Range("S14").Value = 0.05934870938780279
$ws.Range("T14").Value = 0.05934870938780282
$ws.Range("G15").Value = 67.370907
$ws.Range("H15").Value = 202.112721
$ws.Range("I15").Value = 0.2400398653924534
$ws.Range("J15").Value = 0.2400398653924535
$ws.Range("M15").Value = 12.24131666666667
$ws.Range("N15").Value = 36.72395
$ws.Range("O15").Value = 0.2084884078209579
$ws.Range("P15").Value = 0.2084884078209579
$ws.Range("Q15").Value = 824.7086067075501
$ws.Range("R15").Value = 7422.377460367951
$ws.Range("S15").Value = 0.05004552934922966
$ws.Range("T15").Value = 0.05004552934922968
$ws.Range("G16").Value = 67.370907
$ws.Range("H16").Value = 202.112721
$ws.Range("I16").Value = 0.2400398653924534
$ws.Range("J16").Value = 0.2400398653924535
$ws.Range("M16").Value = 20.799674
$ws.Range("N16").Value = 62.399022
$ws.Range("O16").Value = 0.3542503664873991
$ws.Range("P16").Value = 0.3542503664873992
$ws.Range("Q16").Value = 1401.292902684318
$ws.Range("R16").Value = 12611.63612415886
$ws.Range("S16").Value = 0.08503421028686256
$ws.Range("T16").Value = 0.08503421028686262
$ws.Range("G17").Value = 8.783890333333334
$ws.Range("H17").Value = 26.351671
$ws.Range("I17").Value = 0.03129665232554173
$ws.Range("J17").Value = 0.03129665232554174
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 1.222961333333333
$ws.Range("N17").Value = 3.668884
$ws.Range("O17").Value = 0.02082890820948692
$ws.Range("P17").Value = 0.02082890820948692
$ws.Range("Q17").Value = 10.74235823390711
$ws.Range("R17").Value = 96.68122410516401
$ws.Range("S17").Value = 0.000651875098552934
$ws.Range("T17").Value = 0.0006518750985529344
$ws.Range("G18").Value = 8.783890333333334
$ws.Range("H18").Value = 26.351671
$ws.Range("I18").Value = 0.03129665232554173
$ws.Range("J18").Value = 0.03129665232554174
$ws.Range("O18").Value = 0.1691870972318839
$ws.Range("P18").Value = 0.169187097231884
$ws.Range("Q18").Value = 87.25701744616511
$ws.Range("R18").Value = 785.3131570154859
$ws.Range("S18").Value = 0.005294989760033896
$ws.Range("T18").Value = 0.005294989760033899
$ws.Range("G19").Value = 8.783890333333334
$ws.Range("H19").Value = 26.351671
$ws.Range("I19").Value = 0.03129665232554173
$ws.Range("J19").Value = 0.03129665232554174
$ws.Range("M19").Value = 14.516908
$ws.Range("N19").Value = 43.550724
$ws.Range("O19").Value = 0.247245220250272
$ws.Range("P19").Value = 0.2472452202502721
$ws.Range("Q19").Value = 127.5149278510894
$ws.Range("R19").Value = 1147.634350659804
$ws.Range("S19").Value = 0.007737947697324754
$ws.Range("T19").Value = 0.007737947697324757
$ws.Range("G20").Value = 8.783890333333334
$ws.Range("H20").Value = 26.351671
$ws.Range("I20").Value = 0.03129665232554173
$ws.Range("J20").Value = 0.03129665232554174
$ws.Range("M20").Value = 12.24131666666667
$ws.Range("N20").Value = 36.72395
$ws.Range("O20").Value = 0.2084884078209579
$ws.Range("P20").Value = 0.2084884078209579
$ws.Range("Q20").Value = 107.5263831356056
$ws.Range("R20").Value = 967.73744822045
$ws.Range("S20").Value = 0.006524989213478274
$ws.Range("T20").Value = 0.006524989213478277
$ws.Range("G21").Value = 8.783890333333334
$ws.Range("H21").Value = 26.351671
$ws.Range("I21").Value = 0.03129665232554173
$ws.Range("J21").Value = 0.03129665232554174
$ws.Range("M21").Value = 20.799674
$ws.Range("N21").Value = 62.399022
$ws.Range("O21").Value = 0.3542503664873991
$ws.Range("P21").Value = 0.3542503664873992
$ws.Range("Q21").Value = 182.7020553850847
$ws.Range("R21").Value = 1644.318498465762
$ws.Range("S21").Value = 0.01108685055615187
$ws.Range("T21").Value = 0.01108685055615188
$ws.Range("G22").Value = 51.550662
$ws.Range("H22").Value = 154.651986
$ws.Range("I22").Value = 0.1836729609024243
$ws.Range("J22").Value = 0.1836729609024243
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 1.222961333333333
$ws.Range("N22").Value = 3.668884
$ws.Range("O22").Value = 0.02082890820948692
$ws.Range("P22").Value = 0.02082890820948692
$ws.Range("Q22").Value = 63.044466333736
$ws.Range("R22").Value = 567.400197003624
$ws.Range("S22").Value = 0.003825707243201274
$ws.Range("T22").Value = 0.003825707243201276
$ws.Range("G23").Value = 51.550662
$ws.Range("H23").Value = 154.651986
$ws.Range("I23").Value = 0.1836729609024243
$ws.Range("J23").Value = 0.1836729609024243
$ws.Range("O23").Value = 0.1691870972318839
$ws.Range("P23").Value = 0.169187097231884
$ws.Range("Q23").Value = 512.0916635793639
$ws.Range("R23").Value = 4608.824972214275
$ws.Range("S23").Value = 0.03107509509506647
$ws.Range("T23").Value = 0.03107509509506649
$ws.Range("G24").Value = 51.550662
$ws.Range("H24").Value = 154.651986
$ws.Range("I24").Value = 0.1836729609024243
$ws.Range("J24").Value = 0.1836729609024243
$ws.Range("M24").Value = 14.516908
$ws.Range("N24").Value = 43.550724
$ws.Range("O24").Value = 0.247245220250272
$ws.Range("P24").Value = 0.2472452202502721
$ws.Range("Q24").Value = 748.3562175930959
$ws.Range("R24").Value = 6735.205958337864
$ws.Range("S24").Value = 0.04541226167233949
$ws.Range("T24").Value = 0.04541226167233951
$ws.Range("G25").Value = 51.550662
$ws.Range("H25").Value = 154.651986
$ws.Range("I25").Value = 0.1836729609024243
$ws.Range("J25").Value = 0.1836729609024243
$ws.Range("M25").Value = 12.24131666666667
$ws.Range("N25").Value = 36.72395
$ws.Range("O25").Value = 0.2084884078209579
$ws.Range("P25").Value = 0.2084884078209579
$ws.Range("Q25").Value = 631.0479779183
$ws.Range("R25").Value = 5679.4318012647
$ws.Range("S25").Value = 0.03829368317830748
$ws.Range("T25").Value = 0.0382936831783075
$ws.Range("G26").Value = 51.550662
$ws.Range("H26").Value = 154.651986
$ws.Range("I26").Value = 0.1836729609024243
$ws.Range("J26").Value = 0.1836729609024243
$ws.Range("M26").Value = 20.799674
$ws.Range("N26").Value = 62.399022
$ws.Range("O26").Value = 0.3542503664873991
$ws.Range("P26").Value = 0.3542503664873992
$ws.Range("Q26").Value = 1072.236964084188
$ws.Range("R26").Value = 9650.132676757692
$ws.Range("S26").Value = 0.06506621371350953
$ws.Range("T26").Value = 0.06506621371350955
